# "Generate Report for Handoff"
#
# The localization status moved from "In Translation" to "Ready for
# handoff" and the HO/handoff timestamps were refreshed to reflect the
# new report generation time. Apply the same update to the Overview
# rollup sheet and to each per-language detail sheet (zh-cn, de-de), then
# widen the "Status" column(s) so the longer label isn't truncated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) / F (de-de) are the Status cells,
#     column G is the "Latest HO Xliff Generate Date" -----------------------
$ws = $wb.Sheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-31 11:04:44"
$ws.Range("E1:F1").EntireColumn.ColumnWidth = 16.3333333333333

# --- zh-cn detail sheet: column C is Status, column H is the Latest
#     Handoff Datetime ------------------------------------------------------
$ws = $wb.Sheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-31 11:04:40"
$ws.Range("C1").EntireColumn.ColumnWidth = 16.3333333333333

# --- de-de detail sheet: column C is Status, column H is the Latest
#     Handoff Datetime ------------------------------------------------------
$ws = $wb.Sheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-31 11:04:44"
$ws.Range("C1").EntireColumn.ColumnWidth = 16.3333333333333
